# Applies the "codes and speedup updated" change:
#  - Fills in measured C/D values for rows 3-8 (Codes table)
#  - Adds speedup formulas (C/D) for rows 14-19 (Speedup table)
#  - Moves the active selection on the sheet to D20

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Codes table: raw measured values -------------------------------------
$ws.Range("C3").Value = 20530
$ws.Range("D3").Value = 7343.4

$ws.Range("C4").Value = 12323
$ws.Range("D4").Value = 3917

$ws.Range("C5").Value = 86954
$ws.Range("D5").Value = 66318

$ws.Range("C6").Value = 20577
$ws.Range("D6").Value = 7585.1

$ws.Range("C7").Value = 1276.2
$ws.Range("D7").Value = 90.274000000000001

$ws.Range("C8").Value = 349.24079999999998
$ws.Range("D8").Value = 10573.045700000001

# --- Speedup table: formulas referencing the codes table ------------------
$ws.Range("D14").Formula = "=C3/D3"
$ws.Range("D15").Formula = "=C4/D4"
$ws.Range("D16").Formula = "=C5/D5"
$ws.Range("D17").Formula = "=C6/D6"
$ws.Range("D18").Formula = "=C7/D7"
$ws.Range("D19").Formula = "=C8/D8"

# --- Sheet selection -------------------------------------------------------
$ws.Range("D20").Select()
